$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.324.62'
$ws.Range("E2").Value = '  -0.65%  '

$ws.Range("D3").Value = '1.809.79'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = '  -0.39%  '

$ws.Range("D5").Value = "'313.28"
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = '  -0.40%  '

$ws.Range("D7").Value = "'0.5141"
$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("E8").Value = '  +3.08%  '

$ws.Range("D9").Value = "'0.07861"
$ws.Range("E9").Value = '  -5.27%  '

$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").Value = "'40.79"
$ws.Range("E11").Value = '  -2.85%  '

$ws.Range("D12").Value = "'6.387"
$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("D13").Value = "'0.9996"
$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").Value = "'20.36"
$ws.Range("E14").Value = '  -4.09%  '

$ws.Range("D15").Value = "'7.360"
$ws.Range("E15").Value = '  -2.03%  '

$ws.Range("D16").Value = '1.803.61'
$ws.Range("E16").Value = '  -1.29%  '

$ws.Range("D17").Value = "'92.90"
$ws.Range("E17").Value = '  -1.38%  '

$ws.Range("E18").Value = '  -3.58%  '

$ws.Range("D19").Value = "'0.06572"
$ws.Range("E19").Value = '  -1.29%  '

$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("D21").Value = "'17.36"
$ws.Range("E21").Value = '  -2.70%  '

$ws.Range("D22").Value = "'6.039"
$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").Value = '28.383.51'
$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = '  -2.16%  '

$ws.Range("E25").Value = '  -1.42%  '

$ws.Range("D26").Value = "'160.98"
$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("D27").Value = "'20.54"
$ws.Range("E27").Value = '  -2.90%  '

$ws.Range("D28").Value = '2.015.77'
$ws.Range("E28").Value = '  -1.03%  '

$ws.Range("D29").Value = "'2.404"
$ws.Range("E29").Value = '  -0.35%  '

$ws.Range("D30").Value = "'129.01"
$ws.Range("E30").Value = '  +2.27%  '

$ws.Range("E31").Value = '  -0.90%  '

$ws.Range("D32").Value = "'1.063"
$ws.Range("E32").Value = '  -3.47%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'3.662"
$ws.Range("E33").Value = '  -0.63%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'5.589"
$ws.Range("E34").Value = '  -2.67%  '

$ws.Range("D35").Value = "'0.07199"
$ws.Range("E35").Value = '  -5.40%  '

$ws.Range("D36").Value = "'9.170"
$ws.Range("E36").Value = '  +4.25%  '

$ws.Range("E37").Value = '  -1.47%  '

$ws.Range("D38").Value = "'0.2179"
$ws.Range("E38").Value = '  -2.50%  '

$ws.Range("D39").Value = "'5.070"
$ws.Range("E39").Value = '  -4.36%  '

$ws.Range("D40").Value = "'11.60"
$ws.Range("E40").Value = '  -2.51%  '

$ws.Range("D41").Value = "'0.6212"
$ws.Range("E41").Value = '  -2.98%  '

$ws.Range("D42").Value = "'0.9991"
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("E43").Value = '  -2.95%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'13.25"
$ws.Range("E44").Value = '  -2.54%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.6022"
$ws.Range("E45").Value = '  -2.51%  '

$ws.Range("D46").Value = "'1.313"
$ws.Range("E46").Value = '  -5.84%  '

$ws.Range("D47").Value = "'3.741"
$ws.Range("E47").Value = '  -1.65%  '

$ws.Range("D48").Value = "'125.74"
$ws.Range("E48").Value = '  -1.59%  '

$ws.Range("D49").Value = "'1.221"
$ws.Range("E49").Value = '  +1.35%  '

$ws.Range("D50").Value = "'1.936"
$ws.Range("E50").Value = '  -3.29%  '

$ws.Range("D51").Value = "'0.06853"
$ws.Range("E51").Value = '  -1.85%  '
